$d = $word.ActiveDocument
$t = $d.Tables(1)
$values = @(
    '16+21=37',
    '0+33=33',
    '72-7=65',
    '74-62=12',
    '65-27=38',
    '2+77=79',
    '52+36=88',
    '4+85=89',
    '20+57=77',
    '51-17=34',
    '49-49=0',
    '42+29=71',
    '75-9=66',
    '60-3=57',
    '12+36=48',
    '37+57=94',
    '92-46=46',
    '37-4=33',
    '67-59=8',
    '34-29=5',
    '39+5=44',
    '47-5=42',
    '38+56=94',
    '61-37=24',
    '10+0=10',
    '10+57=67',
    '90-10=80',
    '31-10=21',
    '23+8=31',
    '8+60=68',
    '26+29=55',
    '45-17=28',
    '80+11=91',
    '7+9=16',
    '12+79=91',
    '37+34=71',
    '44+17=61',
    '39+24=63',
    '29-21=8',
    '26+73=99',
    '50-18=32',
    '68+9=77',
    '43+20=63',
    '83-13=70',
    '76-26=50',
    '54-5=49',
    '52-36=16',
    '70-2=68',
    '15+81=96',
    '59-15=44',
    '98-60=38',
    '30+68=98',
    '1+57=58',
    '64+28=92',
    '14+12=26',
    '31+17=48',
    '99-59=40',
    '41+43=84',
    '11+54=65',
    '73+11=84',
    '59+10=69',
    '84-34=50',
    '7+90=97',
    '14+1=15',
    '11+31=42',
    '21+13=34',
    '10+26=36',
    '27+41=68',
    '2+46=48',
    '5+80=85',
    '34+10=44',
    '75-48=27',
    '54+44=98',
    '89-35=54',
    '32+44=76',
    '43+32=75',
    '85-52=33',
    '6+46=52',
    '7+4=11',
    '89-54=35',
    '34-4=30',
    '25+24=49',
    '44+39=83',
    '3+82=85',
    '15+55=70',
    '26+30=56',
    '10+57=67',
    '77-59=18',
    '40+50=90',
    '2+2=4',
    '5+51=56',
    '78-24=54',
    '41+8=49',
    '86+1=87',
    '18-9=9',
    '1+42=43',
    '93-68=25',
    '74+19=93',
    '51+21=72',
    '47+10=57'
)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$i]
        $i++
    }
}
Write-Host "Updated cells:" $i
